$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).Style = "Normal"
}

$ws.Range("D2").Value = "34.185.18"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.788.32"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue $ws "D5" "226.43"
$ws.Range("E5").Value = "  +0.48%  "
Set-TextValue $ws "D8" "31.88"
$ws.Range("E8").Value = "  -0.73%  "
Set-TextValue $ws "D9" "0.292"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "2.046.60"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.802.54"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws "D14" "11.03"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "34.151.17"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("E17").Value = "  +0.59%  "
Set-TextValue $ws "D18" "68.26"
$ws.Range("E18").Value = "  +1.18%  "
Set-TextValue $ws "D19" "246.72"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("E21").Value = "  +0.10%  "
Set-TextValue $ws "D22" "10.83"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  +0.03%  "
Set-TextValue $ws "D25" "161.11"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E26").Value = "  +1.29%  "
Set-TextValue $ws "D27" "16.34"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  -0.16%  "
Set-TextValue $ws "D31" "0.0519"
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").Value = "1.442.76"
$ws.Range("E35").Value = "  +3.77%  "
Set-TextValue $ws "D36" "0.647"
$ws.Range("E36").Value = "  -1.81%  "
Set-TextValue $ws "D37" "2.42"
$ws.Range("E37").Value = "  +7.93%  "
$ws.Range("E38").Value = "  +2.83%  "
$ws.Range("E39").Value = "  +0.03%  "
Set-TextValue $ws "D40" "80.48"
$ws.Range("E40").Value = "  +3.42%  "
$ws.Range("E41").Value = "  +0.55%  "
Set-TextValue $ws "D42" "0.922"
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("E43").Value = "  +1.36%  "
Set-TextValue $ws "D44" "13.50"
$ws.Range("E44").Value = "  +1.14%  "
Set-TextValue $ws "D45" "0.0508"
$ws.Range("E45").Value = "  +2.26%  "
Set-TextValue $ws "D46" "6.06"
$ws.Range("E46").Value = "  +3.68%  "
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("E48").Value = "  -6.94%  "
$ws.Range("D49").Value = "1.947.97"
$ws.Range("E49").Value = "  +0.76%  "
Set-TextValue $ws "D50" "105.60"
$ws.Range("E50").Value = "  -2.62%  "
$ws.Range("E51").Value = "  +0.06%  "
